# parallel: <<lab3>>: update table
#
# Row 7 (the "конвеер" delegation-pipeline row for the 500/500 task size)
# was missing the "Ускорение"/"Эффективность" figures for the 1-processor
# block (H7) and had placeholder "-" text in the 2-processor block
# (I7:K7). Replace the placeholder "-" in H7 with the same G-column/2
# "Ускорение" formula used by every other row, and clear I7:J7:K7 back to
# blank (matching the rest of the sheet, where the 2-processor columns
# for this section are left empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$h7 = $ws.Range("H7")
$h7.Formula = "=G7/2"
# Restore the plain "General"/centered-off numeric look used by the rest
# of column H (style was previously the centered text style shared with
# the "-" placeholders in I7:K7).
$h7.HorizontalAlignment = 1
$h7.NumberFormat = "General"

$ws.Range("I7").Value = $null
$ws.Range("J7").Value = $null
$ws.Range("K7").Value = $null

$ws.Range("H7").Select()
